# Atualização de bases das ligas, do dia: 2024-01-31 às 20-02
#
# The underlying "fixture" rows got reshuffled: each group of rows below
# (all sharing the same match date / gameweek block) has its B:AC payload
# (everything except the running index in column A) cyclically rotated by
# one position: new_row[i] = old_row[i+1] (wrapping around the group).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 29  # column AC

$groups = @(
    ,@(322, 323)
    ,@(329, 330, 331)
    ,@(352, 353)
    ,@(362, 363)
    ,@(382, 383)
    ,@(407, 408)
    ,@(430, 431)
    ,@(454, 455, 456)
    ,@(511, 512)
    ,@(516, 517)
    ,@(688, 689)
)

foreach ($group in $groups) {
    $n = $group.Length

    # Snapshot every cell (row x col) in the group before writing anything,
    # so later writes don't clobber values still needed as a source.
    $snapshot = @{}
    for ($i = 0; $i -lt $n; $i++) {
        $row = $group[$i]
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $key = "$i|$col"
            $snapshot[$key] = $ws.Cells.Item($row, $col).Value2
        }
    }

    # new_row[i] <- old_row[i+1] (wrap around)
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcIdx = ($i + 1) % $n
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $key = "$srcIdx|$col"
            $ws.Cells.Item($destRow, $col).Value2 = $snapshot[$key]
        }
    }
}
